{"js": "// Office.js (Word JavaScript API) edit script.\n// Applies the schedule-table text updates described by the diff:\n//   - \"Filesystems; Version Control; Git\" -> \"Filesystems; Markup Languages; Quarto\"\n//   - \"Data Structure & Wrangling\"         -> \"A Field Guide to Data\"\n//   - \"Markup; Data Formats\"               -> \"Version Control; Git; Data Formats\"\n//   - \"Tidy Data\"                          -> \"Wrangling Tidy Data\"\n\nconst body = context.document.body;\n\nconst replacements = [\n  [\"Filesystems; Version Control; Git\", \"Filesystems; Markup Languages; Quarto\"],\n  [\"Data Structure & Wrangling\", \"A Field Guide to Data\"],\n  [\"Markup; Data Formats\", \"Version Control; Git; Data Formats\"],\n  [\"Tidy Data\", \"Wrangling Tidy Data\"],\n];\n\nfor (const [find, replace] of replacements) {\n  const results = body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const result of results.items) {\n    result.insertText(replace, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# Applies the schedule-table text updates described by the diff:\n#   - \"Filesystems; Version Control; Git\" -> \"Filesystems; Markup Languages; Quarto\"\n#   - \"Data Structure & Wrangling\"         -> \"A Field Guide to Data\"\n#   - \"Markup; Data Formats\"               -> \"Version Control; Git; Data Formats\"\n#   - \"Tidy Data\"                          -> \"Wrangling Tidy Data\"\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\nfor ($r = 1; $r -le $t.Rows.Count; $r++) {\n    for ($c = 1; $c -le $t.Columns.Count; $c++) {\n        $cell = $t.Cell($r, $c)\n        # Cell.Range.Text carries a trailing cell-mark (CR + BEL); strip it\n        # before comparing against the plain target strings.\n        $txt = $cell.Range.Text.TrimEnd([char]13, [char]7)\n\n        if ($txt -eq \"Filesystems; Version Control; Git\") {\n            $cell.Range.Text = \"Filesystems; Markup Languages; Quarto\"\n        } elseif ($txt -eq \"Data Structure & Wrangling\") {\n            $cell.Range.Text = \"A Field Guide to Data\"\n        } elseif ($txt -eq \"Markup; Data Formats\") {\n            $cell.Range.Text = \"Version Control; Git; Data Formats\"\n        } elseif ($txt -eq \"Tidy Data\") {\n            $cell.Range.Text = \"Wrangling Tidy Data\"\n        }\n    }\n}\n"}
